$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Rename sheets to the eLife-formatted figure names
$ws1.Name = "growth-restored cell lineage"
$ws2.Name = "growth-halted cell lineage"
$ws3.Name = "non-deleted cell lineage"

# Sheet 1 ("growth-restored cell lineage"): zoom changes from 73 -> 70,
# leave the existing selection (E185) untouched.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 70

# Sheet 3 ("non-deleted cell lineage") becomes the active/selected tab,
# with its selection collapsed from E521:E530 to the single cell E546.
$ws3.Activate()
$ws3.Range("E546").Select()
